# Diary workbook update — add six new diary entries (rows 30-35) covering
# 2020-03-05 through 2020-03-16, matching the final weeks of the quarter
# (testing cases, pull-request work, and class wrap-up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 30 — 2020-03-05
# ---------------------------------------------------------------------
$ws.Range("A30").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A30").Value = 43895
$ws.Range("B30").Value = "classtime"
$ws.Range("C30").Value = "self"
$ws.Range("D30").Value = "Learn more about key expert practice and get feedback regarding pull request "
$ws.Range("E30").Value = "Learn about mechanism for code understanding through reading testing cases"
$ws.Range("F30").Value = "Simply reading test cases can help us understand a lot more about the system the test case is aiming for. For deeper understanding, trying to write test cases ourselves can really help us understand how this feature really works in the system. Apart from reading and writing test cases, for smaller projects, I would use other tactics such as print statement or loggings to see what is really happening in the code."
$ws.Range("G30").Value = "feel good that stuff learned from other class (testing class) can be applied to this class. "
$ws.Rows.Item(30).RowHeight = 171.75

# ---------------------------------------------------------------------
# Row 31 — 2020-03-07
# ---------------------------------------------------------------------
$ws.Range("A31").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A31").Value = 43897
$ws.Range("B31").Value = "1:00 - 4:00"
$ws.Range("C31").Value = "teammates (online)"
$ws.Range("D31").Value = "study existing test cases of Cassandra"
$ws.Range("E31").Value = "Scanned through almost all test cases and found several interesting ones."
$ws.Range("F31").Value = "We feel it was almost impossible to understand all test cases of Cassandra, just like we can not understand all source code. However, by reading through test cases, we are more clear about how the some feature actually works. We gain better understanding of Cassandra by reading testing codes."
$ws.Range("G31").Value = "Always feel good when learned more but still worried about the last pull request."
$ws.Rows.Item(31).RowHeight = 154.5

# ---------------------------------------------------------------------
# Row 32 — 2020-03-12
# ---------------------------------------------------------------------
$ws.Range("A32").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A32").Value = 43902
$ws.Range("B32").Value = "classtime"
$ws.Range("C32").Value = "self"
$ws.Range("D32").Value = "wrap-up of class"
$ws.Range("E32").Value = "reviewed stuff we learned through the quarter especially key expert practice and got a lot of useful advice for my coding career"
$ws.Range("F32").Value = "Before this class I realize that as a programmer we should spend most of my time writing code. However, my perspective has changed after a quarter of this class material and all the guest speaking. I realize that even professionals spend a lot of time reading either their own or other's code for different reasons. I do not try to avoid reading code now but instead I choose to read and learn from the code."
$ws.Range("G32").Value = "Feel really good about how to really become an expert in coding"
$ws.Rows.Item(32).RowHeight = 183.75

# ---------------------------------------------------------------------
# Row 33 — 2020-03-14
# ---------------------------------------------------------------------
$ws.Range("A33").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A33").Value = 43904
$ws.Range("B33").Value = "1:00 - 4:00"
$ws.Range("C33").Value = "teammates (online)"
$ws.Range("D33").Value = "work on last pull request"
$ws.Range("E33").Value = "Our pull request got approved by Kaj, which means we can really work on it now "
$ws.Range("F33").Value = "Finding a pull request that we can actually work on is not too hard (as we only understand part of the whole system, there are limited choices for us)"
$ws.Range("G33").Value = "feel good that our pull request got approved, but also worried if we can actually fix that issue."
$ws.Rows.Item(33).RowHeight = 142.5

# ---------------------------------------------------------------------
# Row 34 — 2020-03-15
# ---------------------------------------------------------------------
$ws.Range("A34").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A34").Value = 43905
$ws.Range("B34").Value = "7:00 pm - 12:00"
$ws.Range("C34").Value = "Tianyu Qi (online)"
$ws.Range("D34").Value = "work on last pull request"
$ws.Range("E34").Value = "finshed part of the fix and there are so much more to do. Wrote a new test case."
$ws.Range("F34").Value = "The seemingly straightforward issue can really involved a bunch of code. While trying to fix the pull request, we realize that there are more than 20 files/methods we need to change. Also from the previous research, it seems like only a few out of hundreds of pull requests got merged, however, this is still a good practice for us even we do not expect developers to merge our changes."
$ws.Range("G34").Value = "feel worried not being able to finish it all. "
$ws.Rows.Item(34).RowHeight = 177

# ---------------------------------------------------------------------
# Row 35 — 2020-03-16
# ---------------------------------------------------------------------
$ws.Range("A35").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("A35").Value = 43906
$ws.Range("B35").Value = "3:00 pm - 9:00 "
$ws.Range("C35").Value = "teammates (online)"
$ws.Range("D35").Value = "work on last pull request"
$ws.Range("E35").Value = "while still working on the last pull request, we  wrote report on our experience for the last homework."
$ws.Range("F35").Value = "While reading the source code, I was happy to see some design pattern that I recently just learned about. Again design pattern is something that is really important when reading people's code or write my own code in an elegant way.`n"
$ws.Range("G35").Value = "feel excited about not having to take the test. Ready to work on the pull request and try to fix issue even after the quarter ends"
$ws.Rows.Item(35).RowHeight = 132.75

# ---------------------------------------------------------------------
# View state — approximate the author's last-saved window position
# (zoom level, scroll position and active selection). These are purely
# cosmetic (not part of the saved data) but are applied on a best-effort
# basis to mirror the commit.
# ---------------------------------------------------------------------
$ws.Range("E35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 85
